# Ex6-RatioAnalysis.xlsx: refreshed Quick Sort comparison counts (column C)
# on the "Quick Sort" sheet; dependent ratio-analysis formulas in columns
# D:M recalc automatically from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quick Sort")
$ws.Activate()

$newComparisons = @{
    4  = 0
    5  = 16
    6  = 194
    7  = 593
    8  = 3914
    9  = 9225
    10 = 60487
    11 = 146965
}

foreach ($row in $newComparisons.Keys) {
    $ws.Cells.Item($row, 3).Value = $newComparisons[$row]
}

# Leave the selection where the author left it when the workbook was saved.
$ws.Range("C12").Select()
